$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value (45188 -> 2023-09-19) for
# every data row (rows 2 through 236). This automatic update bumps that
# value by one day (45188 -> 45189, i.e. 2023-09-20) for every row.
$lastRow = 236
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45189
}
